$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.947.34'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.598.26'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.21'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.594.41'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.127'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.25'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +5.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.395'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.195.45'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.53'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.02%  '
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.585.25'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.939.84'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.16'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.73'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.91'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '398.39'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.594'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.731.14'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.24'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000119'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.21'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.64'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +26.22%  '
$ws.Range('E30').Value = '  +2.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.64'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.591.77'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.66'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.99%  '
$ws.Range('E35').Value = '  +1.53%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.43'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +7.94%  '
$ws.Range('E38').Value = '  +4.81%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.08'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '168.01'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0842'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.841'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '27.14'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.29'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +7.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.15'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.57'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.72'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.04'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.463.23'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.911'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +9.03%  '
